$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with new test case name and new numeric values
$ws.Range("A2").Value = "_MultipleBugs_.NOB_1.ID_165"
$ws.Range("B2").Value = "Empty.PL_Interface_impl.21"
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 257
$ws.Range("F2").Value = 11
$ws.Range("G2").Value = 11
$ws.Range("H2").Value = 440
$ws.Range("I2").Value = 11
$ws.Range("J2").Value = 446
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 40
$ws.Range("M2").Value = 446

# Delete the old row 3 entirely (no longer needed results)
$ws.Rows("3:3").Delete()
